$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Udabda56058ad42c40f27c517b34c3c48"
$ws.Range("B2").Value = "bomb"

$ws.Range("C4").Select()
